# Update the "Login Credentials" sheet: expand the existing last row's
# "D" suffix counter and append a brand-new row of test login data
# (mirrors how the TestNG data-provider rows were previously extended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

# Row 32: bump the trailing " D" markers by one more "D", and append a
# trailing comma to the order-number list.
$ws.Range("D32").Value = "Ahluwalia D D"
$ws.Range("E32").Value = "Sucheta Ahluwalia D D"
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = ",,000051821,000051826,000051828,000051830,"
$ws.Range("F32").Style = "Normal"

# New row 33: a fresh set of login credentials appended to the table.
$ws.Range("A33").Value = "xjeic@gmail.com"
$ws.Range("B33").Value = "U9h49w153@"
$ws.Range("C33").Value = "Nalini"
$ws.Range("D33").Value = "Sethi D"
$ws.Range("E33").Value = "Nalini Sethi D"
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = ",,,000052964"
$ws.Range("F33").Style = "Normal"
